# Generate Report for Archive
#
# Handback processing moved this package from "Ready for handoff" into
# "In Translation", and recorded the first handoff's name
# ("TestHandoff1") against both the zh-cn and de-de localization rows.
# The Overview sheet mirrors each language's Status in its per-locale
# column, so it picks up the same text change.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2 / F2)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

# zh-cn sheet: Status (C2) + Lastest Handoff Name (I2)
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("I2").Value = "TestHandoff1"
$wsZhCn.Columns.Item(3).AutoFit()

# de-de sheet: Status (C2) + Lastest Handoff Name (I2)
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("I2").Value = "TestHandoff1"
$wsDeDe.Columns.Item(3).AutoFit()
